# Sample Project / Main.xlsx - "Project Sample Project is saved.TEST" edit
#
# The rule-table row 11 ("R40" rule) had its Rule-name cell (B11) changed
# from the text "R40" to the text "1". The new literal is still a text
# label (same as the other rule names R10/R20/R30 above it), not a number,
# so it must be written as a shared string, not a numeric value.
#
# Excel's default type inference would turn a bare "1" into a number, so
# the cell is first marked as Text (format code "@") before the value is
# written; this keeps the COM Value assignment from silently re-typing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
